# edit.ps1 -- replays the authored change against before.pptx
#
# 1) Slide 16's table ("Google Shape;213;p29") switches from the deck's
#    custom "Table_0" style to the built-in table style
#    {ACF8778D-2D10-49A1-A3FA-EB6C84009F1C}.
# 2) The presentation's theme ("Integral") is swapped for the default
#    "Office Theme" colour set (the authored diff trades theme1.xml's
#    content for theme2.xml's and vice-versa; theme2.xml is the part
#    actually bound to the slide master, so re-pointing its colour
#    scheme to the stock Office palette reproduces the visible effect).

$p = $ppt.ActivePresentation

# --- 1) table style id -------------------------------------------------
$targetSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $targetSlide.Shapes.Count; $i++) {
    $shp = $targetSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{ACF8778D-2D10-49A1-A3FA-EB6C84009F1C}")
    }
}

# --- 2) theme colours ---------------------------------------------------
function ColorVal($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order matches MsoThemeColorSchemeIndex:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = ColorVal($officeTheme[$i - 1])
}
